$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Cells.Item(1000, 26)

# Column A: sequential match numbers (358-380) as text, added first
$numbers = @("358","359","360","361","362","363","364","365","366","367","368","369","370","371","372","373","374","375","376","377","378","379","380")
for ($i = 0; $i -lt $numbers.Length; $i++) {
    $scratch.Formula = '="' + $numbers[$i] + '"'
    $scratch.Copy()
    $ws.Cells.Item(359 + $i, 1).PasteSpecial(-4163)
}
$scratch.Clear()

# Column B: match URLs, added second
$urls = @(
    "https://fbref.com/en/matches/4e211d29/Fulham-Manchester-City-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/16d5b6aa/Wolverhampton-Wanderers-Crystal-Palace-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/17e9dbbf/West-Ham-United-Luton-Town-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/45c36efd/Tottenham-Hotspur-Burnley-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/482e6ce4/Everton-Sheffield-United-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/6dbd7d72/Newcastle-United-Brighton-and-Hove-Albion-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/e7b6ebb6/Bournemouth-Brentford-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/384547b8/Nottingham-Forest-Chelsea-May-11-2024-Premier-League",
    "https://fbref.com/en/matches/ab9827f2/Manchester-United-Arsenal-May-12-2024-Premier-League",
    "https://fbref.com/en/matches/72c0709d/Aston-Villa-Liverpool-May-13-2024-Premier-League",
    "https://fbref.com/en/matches/2ba17e6d/Tottenham-Hotspur-Manchester-City-May-14-2024-Premier-League",
    "https://fbref.com/en/matches/b42b7907/Brighton-and-Hove-Albion-Chelsea-May-15-2024-Premier-League",
    "https://fbref.com/en/matches/e43e8597/Manchester-United-Newcastle-United-May-15-2024-Premier-League",
    "https://fbref.com/en/matches/0fde9d70/Luton-Town-Fulham-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/273a89b4/Sheffield-United-Tottenham-Hotspur-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/29335211/Manchester-City-West-Ham-United-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/3e33bd98/Arsenal-Everton-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/47ecdc19/Brighton-and-Hove-Albion-Manchester-United-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/546e1a3d/Brentford-Newcastle-United-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/7c034003/Chelsea-Bournemouth-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/c975c7a6/Crystal-Palace-Aston-Villa-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/d4823ed5/Liverpool-Wolverhampton-Wanderers-May-19-2024-Premier-League",
    "https://fbref.com/en/matches/dafb05fe/Burnley-Nottingham-Forest-May-19-2024-Premier-League"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $ws.Cells.Item(359 + $i, 2).Value = $urls[$i]
}
